# Append a new data row (row 86) to Sheet1, following the same pattern as
# the existing rows: Date in column A (as text, not an auto-converted date
# serial), and the three numeric series values in columns B, C, D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86

# Format column A as Text first so the date-like string "2025-11-09" is
# stored as a literal string value instead of being auto-parsed into a
# date serial number. Then reset the style back to Normal so the cell
# doesn't carry an extra/unwanted style index.
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("A" + $row).Value = "2025-11-09"
$ws.Range("A" + $row).Style = "Normal"

$ws.Range("B" + $row).Value = 57.38000106811523
$ws.Range("C" + $row).Value = 405.7000122070312
$ws.Range("D" + $row).Value = 306.1000061035156
